# Update by Rolando's change
#
# - Adds three new TASK rows (45-47) for Dionis covering the
#   "Generar xml de una camapaña" user history plus the follow-up
#   "Genear Apk y Acualizable de una campana" history.
# - Fills in the previously-empty progress cells (C11:F11) on the
#   CONTENT_GENERATOR sheet for the "Generar xml de una camapaña" row.
# - Leaves the workbook with TASK as the active sheet/tab, mirroring the
#   updated selections recorded in the sheet views.

$xlPasteFormats = -4122

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TASK")
$ws2 = $wb.Worksheets.Item("CONTENT_GENERATOR")

# ---------------------------------------------------------------------
# TASK sheet: three new rows at the bottom of the table (45-47)
# ---------------------------------------------------------------------

$taskTexts = @(
  "Implementar CRUD en el Backed de `"Generar xml de una camapaña`"",
  "Integrar servicios del Backed y el FrontEnd para `"Generar xml de una camapaña`"",
  "Crear interfaz para crear `"Generar xml de una camapaña`""
)

$owner   = "Dionis"
$history = "Genear Apk y Acualizable de una campana"
$startDate = 44306
$endDate   = 44314

# Column B: write all three task descriptions first so the new shared
# strings land in the same order as the rest of the row data.
for ($i = 0; $i -lt 3; $i++) {
  $r = 45 + $i
  $ws1.Range("B$r").Value2 = $taskTexts[$i]
}

# Column C: owner text first (shared string must land before the
# "history" string below to match the diff's new shared-strings order);
# red-font styling is applied afterwards so the new date cellXf (below)
# is created before the new red-font cellXf.
for ($i = 0; $i -lt 3; $i++) {
  $r = 45 + $i
  $ws1.Range("C$r").Value2 = $owner
}

# Column D: follow-up user history text.
for ($i = 0; $i -lt 3; $i++) {
  $r = 45 + $i
  $ws1.Range("D$r").Value2 = $history
}

# Columns E/F: dates. Build the bordered/no-fill date style exactly once
# (format-copy the existing border-only style, then apply the date
# number format), then propagate that already-built style to every other
# date cell via a plain format copy so no duplicate cellXfs get created.
# (Done before the Owner column's font colouring below so the new
# cellXfs come out in the same order as the diff: date style first, then
# the red-font style.)
$ws2.Range("A11").Copy() | Out-Null
$ws1.Range("E45").PasteSpecial($xlPasteFormats) | Out-Null
$ws1.Range("E45").NumberFormat = "mm-dd-yy"
$ws1.Range("E45").Value2 = $startDate

foreach ($addr in @("F45", "E46", "F46", "E47", "F47")) {
  $ws1.Range("E45").Copy() | Out-Null
  $ws1.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
}
$ws1.Range("F45").Value2 = $endDate
$ws1.Range("E46").Value2 = $startDate
$ws1.Range("F46").Value2 = $endDate
$ws1.Range("E47").Value2 = $startDate
$ws1.Range("F47").Value2 = $endDate

$excel.CutCopyMode = 0

# Now colour the Owner column with the new red-font cellXf.
for ($i = 0; $i -lt 3; $i++) {
  $r = 45 + $i
  $ws1.Range("C$r").Font.Color = 255
}

# ---------------------------------------------------------------------
# CONTENT_GENERATOR sheet: fill in row 11 (Generar xml de una camapaña)
# ---------------------------------------------------------------------

# C11/D11 reuse the very same date style just created on TASK!E45.
$ws1.Range("E45").Copy() | Out-Null
$ws2.Range("C11").PasteSpecial($xlPasteFormats) | Out-Null
$ws1.Range("E45").Copy() | Out-Null
$ws2.Range("D11").PasteSpecial($xlPasteFormats) | Out-Null
$ws2.Range("C11").Value2 = $startDate
$ws2.Range("D11").Value2 = $endDate

# E11 reuses the existing percent+border style already used by E2:E4.
$ws2.Range("E2").Copy() | Out-Null
$ws2.Range("E11").PasteSpecial($xlPasteFormats) | Out-Null
$ws2.Range("E11").Value2 = 0.2

# F11 reuses the existing fill+border style already used by F8/F9 (which
# also already hold the very same note text).
$ws2.Range("F9").Copy() | Out-Null
$ws2.Range("F11").PasteSpecial($xlPasteFormats) | Out-Null
$ws2.Range("F11").Value2 = "Falta Actions en Api, Crear Servicio y finalizar la UI."

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Selections / active sheet, matching the recorded view state
# ---------------------------------------------------------------------

$ws2.Activate()
$ws2.Range("C11:D11").Select() | Out-Null

$ws1.Activate()
$ws1.Range("C45:C47").Select() | Out-Null

Write-Host "Added Dionis rows 45-47 to TASK and filled CONTENT_GENERATOR row 11"
